$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.9494788646698
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 2.614804983139038
$ws.Range("D1").Value = 1.215245485305786
$ws.Range("E1").Value = 0.8852695226669312
